$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank separator row before "Test_6" (was row 7), shifting the
# remaining rows down by one.
$ws.Rows("7:7").Insert()

# Insert a second blank separator row before "Test_9" (now at row 11 after
# the first insert), shifting the trailing rows down by one more.
$ws.Rows("11:11").Insert()

# Move the active selection to match the author's final cursor position.
$ws.Range("C20").Select()
